$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.534.10"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").Value = "1.877.42"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("D4").Value = "'1.024"
$ws.Range("E4").Value = "  +3.26%  "
$ws.Range("D5").Value = "'318.57"
$ws.Range("E5").Value = "  +2.51%  "
$ws.Range("D6").Value = "'1.023"
$ws.Range("E6").Value = "  +2.73%  "
$ws.Range("D7").Value = "'0.5150"
$ws.Range("E7").Value = "  +1.50%  "
$ws.Range("D8").Value = "'0.3977"
$ws.Range("E8").Value = "  +2.90%  "
$ws.Range("D9").Value = "'0.08390"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").Value = "'1.116"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").Value = "'42.15"
$ws.Range("E11").Value = "  +2.65%  "
$ws.Range("D12").Value = "'6.272"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").Value = "'20.57"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.828.53"
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").Value = "'7.252"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "'1.024"
$ws.Range("E16").Value = "  +3.44%  "
$ws.Range("D17").Value = "'0.00001112"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("D18").Value = "'91.33"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").Value = "'0.06788"
$ws.Range("E19").Value = "  +2.25%  "
$ws.Range("D20").Value = "'17.79"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("E21").Value = "  +2.78%  "
$ws.Range("D22").Value = "'5.982"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "28.524.70"
$ws.Range("E23").Value = "  +2.07%  "
$ws.Range("D24").Value = "'11.19"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").Value = "'2.290"
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("D26").Value = "'162.31"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("D27").Value = "2.030.09"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").Value = "'20.89"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").Value = "'2.374"
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("D30").Value = "'127.80"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("D31").Value = "'0.1054"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").Value = "'5.824"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").Value = "'3.653"
$ws.Range("E34").Value = "  +2.13%  "
$ws.Range("D35").Value = "'0.02434"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").Value = "'0.06513"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").Value = "'0.2191"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'8.945"
$ws.Range("E38").Value = "  -5.57%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.286"
$ws.Range("E39").Value = "  +5.80%  "
$ws.Range("D40").Value = "'1.192"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("D41").Value = "'0.6453"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").Value = "'5.041"
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("D43").Value = "'11.28"
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("D44").Value = "'0.6050"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").Value = "'3.738"
$ws.Range("E46").Value = "  +2.68%  "
$ws.Range("D47").Value = "'1.227"
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("D48").Value = "'1.999"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").Value = "'1.211"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").Value = "'122.36"
$ws.Range("E50").Value = "  +2.14%  "
$ws.Range("D51").Value = "'0.06858"
$ws.Range("E51").Value = "  +0.30%  "
